$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column F - copy the bold/centered header style from E1
$ws.Range("F1").Value = "Trening"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$data = @(
    ,@(45684.59177627315, 529.4, 11.06, 2.009791323116847, "10-15", "Duża Gra")
    ,@(45684.59248460648, 590.6, 12.43, 2.506721598761423, "10-15", "Duża Gra")
    ,@(45684.59383530093, 707.3, 11.47, 2.697062117712838, "10-15", "Duża Gra")
    ,@(45684.59130752315, 488.9, 9.55, 1.866980399404252, "5-10", "Duża Gra")
    ,@(45684.59248113426, 590.3, 9.39, 2.211378642490933, "5-10", "Duża Gra")
    ,@(45684.59383298611, 707.1, 9.39, 2.573243090084621, "5-10", "Duża Gra")
    ,@(45684.60051354166, 1284.3, 13.53, 3.674399614334106, "10-15", "Mała Gra")
    ,@(45684.60216053241, 1426.6, 13.88, 3.657331671033586, "10-15", "Mała Gra")
    ,@(45684.60346030092, 1538.9, 13.7, 2.951944419315885, "10-15", "Mała Gra")
    ,@(45684.59839664352, 1101.4, 9.99, 2.797082100595746, "5-10", "Mała Gra")
    ,@(45684.60051006945, 1284, 8.76, 2.858810535499027, "5-10", "Mała Gra")
    ,@(45684.60215706019, 1426.3, 8.96, 2.693960411208014, "5-10", "Mała Gra")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}

# Touch A2 with the lowercase ISO format first (registers numFmt 164 in the
# style table, matching the source workbook), then apply the real uppercase
# timestamp format to the whole column - this final call wins for every cell.
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2:A13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
